# Update the course schedule topics (column H) on Sheet1.
# Commit message: "added project 11 to menu on the left" — the actual
# change updates the week-by-week "topic" column text for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H32").Value = "Python Basics (constants, variables, comments, strings, print)"
$ws.Range("H33").Value = "Operators and Expressions, intro to functions"
$ws.Range("H34").Value = "Functions, input from user"
$ws.Range("H39").Value = "Control Flow (for loops), Dictionaries"
$ws.Range("H41").Value = "Files and strings"
$ws.Range("H42").Value = "2D lists, nested for loops"
$ws.Range("H38").Value = "Control Flow (for loops), mutability, random"
$ws.Range("H43").Value = "Data Structures (tuples)"

# Sheet view tweaks: zoom level and active selection cell.
$ws.Application.ActiveWindow.Zoom = 95
[void]$ws.Range("H44").Select()
